$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Group" header in F1
$ws.Range("F1").Value = "Group"

# Column F data mirrors column A's geom_code value, but is blank on every
# third (geom_code "C") row -- representing the new "Group" related table
# column, where a null/blank value should be handled as blank.
$groupValues = @("A", "B", "", "A", "B", "", "A", "B", "")

for ($i = 0; $i -lt $groupValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $groupValues[$i]
}

# Set the new column's width to match the target workbook (~19.4 OOXML
# character units).
$ws.Columns.Item(6).ColumnWidth = 18.5

# Update the active cell selection to F10, matching the authored workbook.
$ws.Range("F10").Select() | Out-Null
